# Update cryptos list (Price / Volume(1h)) with the latest scrape values.
# Leading apostrophe forces numeric-looking Price strings (e.g. "1.012", "30.642.35")
# to be stored as text, matching the source data which is not a real number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.642.35"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "'2.127.90"
$ws.Range("E3").Value = "  +1.09%  "
$ws.Range("D4").Value = "'1.012"
$ws.Range("E4").Value = "  +0.69%  "
$ws.Range("D5").Value = "'352.73"
$ws.Range("E5").Value = "  +5.50%  "
$ws.Range("E6").Value = "  +0.72%  "
$ws.Range("D7").Value = "'0.5281"
$ws.Range("E7").Value = "  +1.06%  "
$ws.Range("D8").Value = "'0.4549"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'54.09"
$ws.Range("E9").Value = "  +1.55%  "
$ws.Range("D10").Value = "'0.09095"
$ws.Range("E10").Value = "  +1.84%  "
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("D12").Value = "'24.67"
$ws.Range("E12").Value = "  +1.75%  "
$ws.Range("D13").Value = "'2.134.34"
$ws.Range("E13").Value = "  +1.38%  "
$ws.Range("D14").Value = "'6.874"
$ws.Range("E14").Value = "  +0.40%  "
$ws.Range("D15").Value = "'8.127"
$ws.Range("E15").Value = "  +0.69%  "
$ws.Range("D16").Value = "'102.37"
$ws.Range("E16").Value = "  +5.98%  "
$ws.Range("D17").Value = "'0.00001176"
$ws.Range("E17").Value = "  +3.06%  "
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("D19").Value = "'0.06719"
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("D20").Value = "'19.51"
$ws.Range("E20").Value = "  +1.47%  "
$ws.Range("D21").Value = "'1.010"
$ws.Range("E21").Value = "  +0.67%  "
$ws.Range("D22").Value = "'6.368"
$ws.Range("E22").Value = "  +0.46%  "
$ws.Range("D23").Value = "'30.729.38"
$ws.Range("E23").Value = "  +0.78%  "
$ws.Range("D24").Value = "'12.89"
$ws.Range("E24").Value = "  +3.14%  "
$ws.Range("D25").Value = "'2.383"
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("D26").Value = "'2.382.69"
$ws.Range("E26").Value = "  +1.23%  "
$ws.Range("D27").Value = "'22.54"
$ws.Range("E27").Value = "  +1.17%  "
$ws.Range("D28").Value = "'2.566"
$ws.Range("E28").Value = "  +1.18%  "
$ws.Range("D29").Value = "'164.73"
$ws.Range("E29").Value = "  +1.25%  "
$ws.Range("D30").Value = "'136.43"
$ws.Range("E30").Value = "  +2.48%  "
$ws.Range("E31").Value = "  -0.79%  "
$ws.Range("D32").Value = "'0.1085"
$ws.Range("E32").Value = "  +1.09%  "
$ws.Range("D33").Value = "'1.675"
$ws.Range("E33").Value = "  +0.50%  "
$ws.Range("D34").Value = "'6.391"
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("E35").Value = "  +2.06%  "
$ws.Range("D36").Value = "'6.170"
$ws.Range("E36").Value = "  +6.50%  "
$ws.Range("E37").Value = "  -1.38%  "
$ws.Range("D38").Value = "'0.02656"
$ws.Range("E38").Value = "  +2.52%  "
$ws.Range("D39").Value = "'0.06897"
$ws.Range("E39").Value = "  +0.77%  "
$ws.Range("D40").Value = "'0.2326"
$ws.Range("E40").Value = "  +1.09%  "
$ws.Range("D41").Value = "'12.57"
$ws.Range("E41").Value = "  -1.10%  "
$ws.Range("D42").Value = "'0.6924"
$ws.Range("E42").Value = "  +0.51%  "
$ws.Range("D43").Value = "'1.277"
$ws.Range("E43").Value = "  +2.23%  "
$ws.Range("E44").Value = "  +5.03%  "
$ws.Range("E45").Value = "  +1.06%  "
$ws.Range("D46").Value = "'0.6463"
$ws.Range("E46").Value = "  +1.35%  "
$ws.Range("D47").Value = "'3.761"
$ws.Range("E47").Value = "  +2.71%  "
$ws.Range("D48").Value = "'0.00000000367"
$ws.Range("E48").Value = "  +5.15%  "
$ws.Range("D49").Value = "'1.256"
$ws.Range("E49").Value = "  +0.44%  "
$ws.Range("D50").Value = "'0.3478"
$ws.Range("E50").Value = "  +2.58%  "
$ws.Range("E51").Value = "  -0.37%  "
